# "Generate Report for Handoff" — refresh the handoff status / timestamps for
# the two files (ca8dbb89-..., f18d2794-...) that are now ready to be handed
# off again, and record the "stale handback" detail message on the per-locale
# sheets.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"

$errCa8 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5231412db5d77f074497441797b3b019e14018d6/e2e/ca8dbb89-398c-47ed-95ea-42ceffe6eb08.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/27e41090acba070d2db1fdee7844cd1ec18d566b/e2e/ca8dbb89-398c-47ed-95ea-42ceffe6eb08.md."
$errF18 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5231412db5d77f074497441797b3b019e14018d6/e2e/f18d2794-c2e9-4f3d-af90-6f9fb437b93a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/27e41090acba070d2db1fdee7844cd1ec18d566b/e2e/f18d2794-c2e9-4f3d-af90-6f9fb437b93a.md."

# ---- Overview sheet ----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E4").Value = $readyStatus
$overview.Range("F4").Value = $readyStatus
$overview.Range("G4").Value = "2016-11-14 07:25:38"

$overview.Range("E5").Value = $readyStatus
$overview.Range("F5").Value = $readyStatus
$overview.Range("G5").Value = "2016-11-14 07:25:38"

# ---- zh-cn sheet ---------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C4").Value = $readyStatus
$zhcn.Range("H4").Value = "2016-11-14 07:25:20"
$zhcn.Range("P4").Value = $errCa8

$zhcn.Range("C5").Value = $readyStatus
$zhcn.Range("H5").Value = "2016-11-14 07:25:20"
$zhcn.Range("P5").Value = $errF18

$zhcn.Columns.Item(16).ColumnWidth = 39.1665

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C4").Value = $readyStatus
$dede.Range("H4").Value = "2016-11-14 07:25:38"
$dede.Range("P4").Value = $errCa8

$dede.Range("C5").Value = $readyStatus
$dede.Range("H5").Value = "2016-11-14 07:25:38"
$dede.Range("P5").Value = $errF18

$dede.Columns.Item(16).ColumnWidth = 39.1665
